$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '71.737.28'
$ws.Range('E2').Value = '  +3.65%  '
$ws.Range('D3').Value = '3.693.78'
$ws.Range('E3').Value = '  +8.16%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '590.61'
$ws.Range('E5').Value = '  +1.63%  '
$ws.Range('E6').Value = '  +1.64%  '
$ws.Range('D7').Value = '3.683.32'
$ws.Range('E7').Value = '  +8.03%  '
$ws.Range('E8').Value = '  +4.34%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('E10').Value = '  +2.71%  '
$ws.Range('E11').Value = '  +5.16%  '
$ws.Range('D12').Value = '50.07'
$ws.Range('E12').Value = '  +3.84%  '
$ws.Range('D13').Value = '0.0000289'
$ws.Range('E13').Value = '  +2.72%  '
$ws.Range('D14').Value = '4.289.91'
$ws.Range('E14').Value = '  +8.20%  '
$ws.Range('D15').Value = '685.18'
$ws.Range('E15').Value = '  +0.96%  '
$ws.Range('D16').Value = '9.07'
$ws.Range('E16').Value = '  +5.04%  '
$ws.Range('D17').Value = '3.692.65'
$ws.Range('E17').Value = '  +8.21%  '
$ws.Range('D18').Value = '71.826.07'
$ws.Range('E18').Value = '  +3.61%  '
$ws.Range('E19').Value = '  +2.30%  '
$ws.Range('D20').Value = '18.19'
$ws.Range('E20').Value = '  +2.47%  '
$ws.Range('E21').Value = '  +3.56%  '
$ws.Range('D22').Value = '6.42'
$ws.Range('E22').Value = '  +19.80%  '
$ws.Range('E23').Value = '  +4.06%  '
$ws.Range('D24').Value = '17.89'
$ws.Range('E24').Value = '  +5.25%  '
$ws.Range('D25').Value = '103.99'
$ws.Range('E25').Value = '  +3.17%  '
$ws.Range('E26').Value = '  +4.04%  '
$ws.Range('D27').Value = '2.85'
$ws.Range('E27').Value = '  +5.82%  '
$ws.Range('D28').Value = '10.26'
$ws.Range('E28').Value = '  +6.45%  '
$ws.Range('D29').Value = '35.41'
$ws.Range('E29').Value = '  +5.72%  '
$ws.Range('E30').Value = '  +6.33%  '
$ws.Range('E31').Value = '  +7.86%  '
$ws.Range('D32').Value = '4.30'
$ws.Range('E32').Value = '  +15.87%  '
$ws.Range('E33').Value = '  +2.95%  '
$ws.Range('D34').Value = '566.46'
$ws.Range('E34').Value = '  +2.41%  '
$ws.Range('E35').Value = '  +4.48%  '
$ws.Range('D36').Value = '59.55'
$ws.Range('E36').Value = '  +2.62%  '
$ws.Range('D37').Value = '3.758.33'
$ws.Range('E37').Value = '  +4.22%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('E39').Value = '  +3.86%  '
$ws.Range('E40').Value = '  +6.24%  '
$ws.Range('D41').Value = '35.76'
$ws.Range('E41').Value = '  +2.29%  '
$ws.Range('E42').Value = '  +6.33%  '
$ws.Range('E44').Value = '  +5.15%  '
$ws.Range('E45').Value = '  +5.90%  '
$ws.Range('D46').Value = '2.90'
$ws.Range('E46').Value = '  +9.24%  '
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('E48').Value = '  +4.08%  '
$ws.Range('E49').Value = '  +3.39%  '
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').Value = '135.58'
$ws.Range('E51').Value = '  +3.61%  '
